$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has two debt records (rows 2 & 3, item numbers 13 & 14).
# A third record (item number 12) needs to be inserted above them, becoming
# the new row 2, and pushing the existing rows down to rows 3 and 4.
#
# Rather than using Rows.Insert() (which copies formatting and introduces a
# brand-new cell style into styles.xml that the target file does not have),
# we rewrite the three data rows' values/text directly, re-using the
# existing style pattern (col A = style 1, cols B:G = style 2) that already
# covers rows 2 and 3.

# Row 3 <- old row 3 content (item 14), moved down to row 4
$ws.Range("A4").Value = 14
$ws.Range("B4").Value = "—般貸款"
$ws.Range("C4").Value = "—般貸款"
$ws.Range("D4").Value = "陳姿伶"
$ws.Range("E4").Value = "台東縣都蘭農會臺東縣東河鄉都蘭村都蘭"
$ws.Range("F4").Value = 330750
$ws.Range("G4").Value = "96年03月13曰"

# Row 2 -> old row 2 content (item 13), moved down to row 3
$ws.Range("A3").Value = 13
$ws.Range("B3").Value = "貸款"
$ws.Range("C3").Value = "信貸"
$ws.Range("D3").Value = "廖國棟"
$ws.Range("E3").Value = "土地銀行台東分行臺東縣台東市中華路"
$ws.Range("F3").Value = "'1848023"
$ws.Range("G3").Value = "100年07月01曰"

# New row 2 (item 12) - property building record, mirrors the row-1 header values
$ws.Range("A2").Value = 12
$ws.Range("B2").Value = "車貸"
$ws.Range("C2").Value = "車貸"
$ws.Range("D2").Value = "陳姿蛉"
$ws.Range("E2").Value = "台新銀行三重分行新北市三重區正義北路"
$ws.Range("F2").Value = 305434
$ws.Range("G2").Value = "99年08月27日"
